# InitProperty.xlsx - add a new "Force" boolean-flag row to the Property1
# sheet's header block.
#
# The header block (rows 1-9 before this edit) lists one flag-row per
# property of a field definition: Type / Public / Private / Save / Cache /
# Ref / Upload / Desc. This change inserts a new "Force" flag row right
# after "Ref" (i.e. it becomes the new row 8), pushing "Upload", "Desc"
# and every data row below it down by one. All the data rows (formerly
# 10-129, now 11-130) keep their original content - they just shift down.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a brand new row at 8 - this shifts rows 8..129 down to 9..130,
# including the dataValidation ranges (A7:A8 -> A7:A9, B7:J8 -> B7:J9) and
# the sheet dimension (A1:F129 -> A1:F130), all handled automatically by
# the row insert.
$ws.Rows.Item(8).Insert()

# Populate the freshly inserted row 8 with the new "Force" flag, matching
# the pattern used by the existing flag rows (label in column A, then a
# FALSE boolean in each of columns B through F).
$ws.Range("A8").Value = "Force"
$ws.Range("B8:F8").Value = $false
